$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "04/08/2023 15:01"
$ws.Range("B2").Value = "Este desayuno con huevos, típico de Oriente Medio, es la cena veraniega perfecta"
$ws.Range("C2").Value = "Miguel Ayuso Rejas"
$ws.Range("F2").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/este-desayuno-huevos-tipico-oriente-medio-cena-veraniega-perfecta"

$ws.Range("A3").Value = "03/08/2023 17:01"
$ws.Range("B3").Value = "Esta comida callejera italiana, sabrosa y fácil, es mi receta favorita con calabacín"
$ws.Range("C3").Value = "Miguel Ayuso Rejas"
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/esta-comida-callejera-italiana-sabrosa-facil-mi-receta-favorita-calabacin"

$ws.Range("A4").Value = "03/08/2023 16:01"
$ws.Range("B4").Value = "Muffins salados de boniato, una receta para chuparse los dedos"
$ws.Range("C4").Value = "Inés Vazquez Noya"
$ws.Range("D4").Value = "Harina integral 200 gLevadura química cucharadita1 Sal y pimienta al gusto Boniato cocido161 gLeche o bebida vegetal de almendras, de avena o de soja300 mlAceite de girasol 30 mlEspinaca fresca 50 gGuisantes 30 gComino molido 1 cuchara aproximadamente10 g"
$ws.Range("E4").Value = "Encender el horno a 190 ºC.En un cuenco mezclar los ingredientes secos: harina integral, levadura química, sal y pimienta. Integrar con un batidor de varilla o un tenedor. Reservar. Por otra parte, en un robot de cocina o licuadora, triturar 100 g de boniato ya cocido con la leche o bebida y aceite de girasol. Pulsar por varios segundos hasta conseguir una bebida liquida de color naranja. Agregar este último al recipiente con la harina e integrar hasta tener una masa homogénea. Momento de agregar la espinaca fresca, los guisantes y los 80 g de boniato restantes en cubitos, más el comino molido. Revolver nuevamente.  Acomodar en molde de muffins de silicona y verter la mezcla anterior sin alcanzar la superficie. De utilizar moldes de acero, pintar con aceite toda la superficie. Cocinar por 30 minutos. Dejar reposar unos 10 minutos antes de desmoldar y disfrutar."
$ws.Range("F4").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/muffins-salados-boniato-receta-para-chuparse-dedos"

$ws.Range("A5").Value = "02/08/2023 11:01"
$ws.Range("B5").Value = "Refrescante, saciante y sin encender un fuego: esta receta de calabacín es una cena saludable de verano perfecta"
$ws.Range("C5").Value = "Liliana Fuchs"
$ws.Range("F5").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/refrescante-saciante-encender-fuego-esta-receta-calabacin-cena-saludable-verano-perfecta"

$ws.Range("A6").Value = "02/08/2023 07:04"
$ws.Range("B6").Value = "Así son los nuevos postres veganos de Dia (que no son yogures) con tres bases vegetales diferentes"
$ws.Range("C6").Value = "Liliana Fuchs"
$ws.Range("F6").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/asi-nuevos-postres-veganos-dia-que-no-yogures-tres-bases-vegetales-diferentes"

$ws.Range("A7").Value = "31/07/2023 09:00"
$ws.Range("B7").Value = "Despedimos el séptimo mes del año en el menú semanal vegetariano del 31 de julio"
$ws.Range("C7").Value = "Inés Vazquez Noya"
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("F7").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/despedimos-septimo-mes-ano-menu-semanal-vegetariano-31-julio"

$ws.Range("A8").Value = "30/07/2023 16:00"
$ws.Range("B8").Value = "Palitos de queso, la receta más fácil y rápida del aperitivo del verano"
$ws.Range("C8").Value = "Carmen Tía Alia"
$ws.Range("D8").Value = "Queso halloumi 200 gMiel 15 mlSemillas de sésamo 5 g"
$ws.Range("E8").Value = "Cortamos el queso halloumi en diez piezas de un tamaño aproximado y con forma de lingote. Los ponemos en remojo con abundante agua y los dejamos reposar en la nevera 30-40 minutos. Con esto conseguiremos rebajar un poco el punto de sal que tiene este queso. Secamos los palitos de halloumi con papel absorbente y los ensartamos en brochetas, aunque también los podemos dejar tal cual. Los marcamos en una parrilla bien caliente por todas sus caras. Después los pincelamos con miel y los espolvoreamos con semillas de sésamo. Servimos inmediatamente."
$ws.Range("F8").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/palitos-queso-receta-facil-rapida-aperitivo-verano"

$ws.Range("A9").Value = "25/07/2023 16:00"
$ws.Range("B9").Value = "Galletas de naranja con harina de garbanzos y avena, una receta fácil para acompañar la merienda"
$ws.Range("C9").Value = "Inés Vazquez Noya"
$ws.Range("D9").Value = "Harina de garbanzo 60 gAvena molida30 gHarina integral 30 gCoco rallado 20 gMantequilla de cacahuete natural 20 mlNaranja zumo y ralladura1 Azúcar moreno 10 gEsencia de vainilla cucharadita1 Agua cantidad necesariaAceite de girasol 20 mlSal pizcaCanela molida 1 cucharadita (opcional)"
$ws.Range("E9").Value = "Encender el horno a 180 ºC.En un cuenco, mezclar los ingredientes secos: la harina de garbanzos, la avena molida, la harina integral, el coco rallado, el azúcar, la levadura química y la ralladura de naranja. Integrar. Agregar luego los líquidos: el aceite de girasol, el zumo de naranja, la esencia de vainilla, la mantequilla de cacahuete y el agua. Sumar una pizca de sal para que realce el sabor. Con ayuda de una cuchara de madera mezclar y terminar de integrar todo con la mano.Con las manos o una cucharada, formar 8 a 10 bolas con la masa y acomodar en una bandeja para el horno con papel de cocina o silpat. De acuerdo al tamaño, será la cantidad de galletas. Aplanar cada bola con un tenedor untado en harina o levemente húmedo. Cocinar por 10 minutos hasta que estén ligeramente doradas en la base. Dejar las galletas en la bandeja durante un minuto y transferir a una rejilla para enfriar antes de servir, así tendrán más cuerpo."
$ws.Range("F9").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/galletas-naranja-harina-garbanzos-avena-receta-facil-para-acompanar-merienda"

$ws.Range("A10").Value = "24/07/2023 10:12"
$ws.Range("B10").Value = "Receta de ensalada coreana de patata o gamja, tan cremosa y refrescante que resulta adictiva"
$ws.Range("C10").Value = "Carmen Tía Alia"
$ws.Range("D10").Value = "Patata mediana2 Zanahoria 1 Pepino 1 Cebolla morada 1 Huevo 2 Maíz en grano (opcional)50 gMayonesa 3 cucharadasVinagre de arroz 2 cucharadasAzúcar 1 cucharadaMostaza de Dijon 1 cucharadita"
$ws.Range("E10").Value = "Pelamos las patatas y las cortamos en trozos de igual tamaño. Calentamos abundante agua con sal en una cacerola y, cuando arranque a hervir, los introducimos junto con los huevos. Retiramos estos últimos después de diez minutos y continuamos cociendo la patata hasta que esté tierna. Pelamos el pepino, la zanahoria y la cebolla morada y picamos todo finamente. Introducimos en un bol y añadimos un poco de sal, removemos y dejamos que se reblandezcan mientras sueltan sus jugos. Pelamos también los huevos y los picamos muy fino.Preparamos la salsa mezclando la mayonesa con el vinagre de arroz, el azúcar y la mostaza. Cuando la patata esté tierna la escurrimos y la pasamos por un pasapurés. Agregamos los huevos duros, las verduras -bien escurridas- y el maíz (opcional).Incorporamos la salsa y removemos para integrar todo bien. Probamos el punto de sal y añadimos más si lo consideramos necesario. Dejamos reposar la ensalada en la nevera, donde la patata irá absorbiendo la salsa poco a poco. Cuando esté bien fría, ya está lista para servir."
$ws.Range("F10").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/receta-ensalada-coreana-patata-gamja-cremosa-refrescante-que-resulta-adictiva"

$ws.Range("A11").Value = "24/07/2023 09:00"
$ws.Range("B11").Value = "Recetas ligerísimas para disfrutar el verano en el menú semanal vegetariano del 24 de julio"
$ws.Range("C11").Value = "Inés Vazquez Noya"
$ws.Range("F11").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/recetas-ligerisimas-para-disfrutar-verano-menu-semanal-vegetariano-24-julio"

$ws.Range("A12").Value = "24/07/2023 07:00"
$ws.Range("B12").Value = "No sabía que mis ensaladas de legumbres tenían un fallo. Todo ha cambiado desde que no me salto este truco"
$ws.Range("C12").Value = "Liliana Fuchs"
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("F12").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/no-sabia-que-mis-ensaladas-legumbres-tenian-fallo-todo-ha-cambiado-que-no-me-salto-este-truco"

$ws.Range("A13").Value = "22/07/2023 10:00"
$ws.Range("B13").Value = "Cómo hacer fatteh, la ensalada de garbanzos especiados, yogur y pan de pita más crujiente de Oriente Medio"
$ws.Range("C13").Value = "Carmen Tía Alia"
$ws.Range("D13").Value = "Pan de pita 3 Garbanzos en conserva, 1 bote grandeYogur griego 250 mlTahini 2 cucharadasZumo de limón 2 Diente de ajo 3 Piñones 25 gComino molido Sal Pimienta negra molida Perejil fresco"
$ws.Range("E13").Value = "Colocamos los garbanzos en una cazuela junto con el líquido de la conserva y un poco de agua, que queden cubiertos a ras. Añadimos el comino, un diente de ajo rallado, sal y pimienta negra. Cocemos 15 minutos o hasta que estén secos. Añadimos un buen puñado de perejil fresco picado. Reservamos. Mientras tanto, cortamos los panes de pita por la mitad, longitudinalmente, y después en triángulos. Los colocamos en una fuente de horno y los tostamos bajo el grill durante unos seis o siete minutos, lo necesario para que queden dorados y bien secos. Reservamos.Por otro lado mezclamos el yogur griego con el tahini, dos dientes de ajo finamente picados (o rallados), el zumo de los dos limones y sal al gusto. Tiene que quedar cremosa, por lo que es posible que necesitemos añadir un poco de agua para aligerarla.Armamos el fatteh colocando los triángulos de pan de pita en la base de una fuente o plato. Extendemos una capa de yogur sobre ellos y luego una capa de garbanzos. Repetimos la operación de nuevo. Espolvoreamos con perejil picado u otra hierba (la menta le va genial), regamos con un chorrito de aceite de oliva virgen extra y servimos."
$ws.Range("F13").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/como-hacer-fatteh-ensalada-garbanzos-especiados-yogur-pan-pita-crujiente-oriente-medio"

$ws.Range("A14").Value = "21/07/2023 16:01"
$ws.Range("B14").Value = "Frittata de verduras en molde de muffins, un clásico italiano para disfrutar en formato individual"
$ws.Range("C14").Value = "Inés Vazquez Noya"
$ws.Range("D14").Value = "Harina de garbanzo 150 gAgua 200 mlBicarbonato sódico 5 gCebolla en polvo 5 gOrégano seco 10 gSal y pimienta al gusto Espinaca fresca 50 gPimiento verde 30 gGuisantes 40 gMaíz amarillo40 g"
$ws.Range("E14").Value = "Encender el horno a 200 ºC.En un cuenco mediano, mezclar la harina de garbanzos junto con la cebolla en polvo, orégano seco, bicarbonato sódico, sal y pimienta. Agregar el agua y batir hasta mezclar. La masa quedará líquida, no espesa. Mientras tanto preparar las verduras, que irán en crudo. El pimiento verde, por ejemplo, en pequeños cubos y las hojas de espinaca groseramente picadas. Junto con el maíz amarillo y los guisantes, añadir a la mezcla anterior e integrar. Con una taza medidora de 1/4, acomodar la masa y completar cada molde de muffin de silicona. De utilizar moldes de acero, pintar con aceite la superficie antes. Hornear durante 30 minutos, hasta que el borde se vea dorado y el centro esté compacto y cuajado. Pasado ese tiempo, apagar el horno y mantener los muffins dentro por 10 minutos más. Una vez fuera, reposar unos minutos y aflojar los lados con una para quitar fácilmente."
$ws.Range("F14").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/frittata-verduras-molde-muffins-clasico-italiano-para-disfrutar-formato-individual"

$ws.Range("A15").Value = "20/07/2023 16:01"
$ws.Range("B15").Value = "Cuscús a la naranja en 10 minutos: una receta con un truco sencillo para darle más sabor al cereal"
$ws.Range("C15").Value = "Inés Vazquez Noya"
$ws.Range("D15").Value = "Cuscús 200 gZumo de naranja 200 mlSal y pimienta al gusto"
$ws.Range("E15").Value = "En un cuenco, volcar el cuscús y agregar el zumo de naranja (levemente templado en el microondas unos dos minutos antes). Tapar y dejar reposar 10 minutos.Con un tenedor, soltar los granos de cuscús y condimentar con la sal y la pimienta al gusto antes de servir."
$ws.Range("F15").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/cuscus-a-naranja-10-minutos-receta-truco-sencillo-para-darle-sabor-al-cereal"

$ws.Range("A16").Value = "19/07/2023 16:01"
$ws.Range("B16").Value = "Trufas de granola, un snack fácil de hacer, nutritivo y bien rico con solo cuatro ingredientes"
$ws.Range("C16").Value = "Inés Vazquez Noya"
$ws.Range("D16").Value = "Granola 150 gMantequilla de cacahuete natural o de almendras - aproxiadamente100 gEsencia de vainilla 10 gChocolate negro o blanco - opcional"
$ws.Range("E16").Value = "En un robot de cocina, procesar la granola unos segundos, hasta obtener un arenado. Volcar en un cuenco, y agregar la esencia de vainilla junto a la mantequilla de cacahuete hasta lograr una masa maleable. Con las manos, tomar pociones de la masa y darle una forma redonda. Acomodar en un plato y de manera opcional: derretir el chocolate negro a baño maría o al microondas, y decorar cada bocado o bañar la trufa entera.  Conservar en la nevera hasta el momento de consumir."
$ws.Range("F16").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/trufas-granola-snack-facil-hacer-nutritivo-bien-rico"

$ws.Range("A17").Value = "19/07/2023 10:01"
$ws.Range("B17").Value = "Ensalada de patatas, pepino y maíz dulce: una combinación de sabores bien rica e ideal para el calor"
$ws.Range("C17").Value = "Inés Vazquez Noya"
$ws.Range("D17").Value = "Patata cocidas y frías200 gPepino 100 gMaíz dulce 80 gPerejil o cilantro o eneldo -- puñado generosoYogur natural o mayonesa80 mlLimón zumo y ralladura1 Aceite de oliva virgen extra Sal y pimienta al gusto"
$ws.Range("E17").Value = "Comenzar hirviendo las patatas enteras en un cazo con abundante agua, hasta que estén tiernas. Una vez listas, escurrir, reservar y dejar enfriar. Con ayuda de un cuchillo, quitar la piel y cortar en láminas, cuadrados grandes o cubos pequeños. Tomar el pepino, cortar a lo largo y luego en mitades. Reservar. En una ensaladera o cuenco profundo, volcar las patatas, el pepino y el maíz. Condimentar con la sal y la pimienta. Agregar el yogur natural y aligerar con el zumo de limón y aceite de oliva virgen extra para lograr una simil salsa o vinagreta. Por último, añadir el perejil picado, la ralladura de limón y revolver una vez más para integrar los sabores. Conservar en la nevera hasta el momento de servir para una mayor frescura."
$ws.Range("F17").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/ensalada-patatas-pepino-maiz-dulce-combinacion-sabores-bien-rica-e-ideal-para-calor"

$ws.Range("A18").Value = "19/07/2023 07:01"
$ws.Range("B18").Value = "Siete trucos para exprimir más zumo de limón más fácilmente sin usar el exprimidor"
$ws.Range("C18").Value = "Liliana Fuchs"
$ws.Range("D18").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("F18").Value = "https://www.directoalpaladar.com/ingredientes-y-alimentos/siete-trucos-para-exprimir-zumo-limon-facilmente-usar-exprimidor"

$ws.Range("A19").Value = "18/07/2023 10:25"
$ws.Range("B19").Value = "Croquetas de arroz integral, zanahoria y curry, una receta que te sacará de más de un apuro"
$ws.Range("C19").Value = "Inés Vazquez Noya"
$ws.Range("D19").Value = "Arroz de grano largo o redondo o integral350 gZanahoria rallada200 gCurry molido 20 gAjo en polvo o cebolla en polvo10 gSal y pimienta al gusto Harina de trigo o de arroz, garbanzos o integral30 gPan rallado o panko para rebozarOrégano seco 15 g"
$ws.Range("E19").Value = "En una olla, cocer el arroz por 20 a 25 minutos, siguiendo las instrucciones del paquete. Reservar hasta enfriar.Tomar la zanahoria y rallar de lado más fino. En un robot de cocina, batidora o procesador de alimentos, colocar el arroz cocido con la zanahoria y triturar unos segundos para que se rompan los granos, sin desintegrar los ingredientes. La intención es que tenga una textura imperfecta y granulada. En un cuenco, mezclar la mezcla anterior y condimentar con el curry, el orégano seco, el ajo, la sal y la pimienta en polvo. Integrar hasta tener una masa homogénea. Incorporar las cucharadas de harina (a elección de cada uno) para absorber cualquier resto de humedad. Al final, la masa se debe poder manejar con las manos. De no ser el caso, incorporar una cucharada de harina, o las necesarias, hasta alcanzar la consistencia optima.Formar las croquetas de 15 a 16 unidades con las manos, dándole la forma deseadas. Intentar que todas tengan el mismo tamaño y dimensión para lograr una cocción uniforme. Pasar por pan rallado o panko para un acabado crujiente. Hornear a 180 ºC con calor arriba y abajo por 30 minutos aproximadamente, hasta que estén doradas."
$ws.Range("F19").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/croquetas-arroz-integral-zanahoria-curry-receta-que-te-sacara-apuro"

$ws.Range("A20").Value = "17/07/2023 10:01"
$ws.Range("B20").Value = "Ensalada de lechuga, maíz dulce y manzana, una receta de verano de Karlos Arguiñano"
$ws.Range("C20").Value = "Inés Vazquez Noya"
$ws.Range("D20").Value = "Hojas verdes 100 gMaíz dulce 50 gManzana roja o verde1 Miel 1 cucharadaVinagre de vino blanco 1 cucharadaAceite de oliva virgen extra 2 a 3 cucharadasSal y pimienta al gusto Orégano seco opcionalSemillas de calabaza opcional - para la decoración"
$ws.Range("E20").Value = "Lavar las hojas verdes y colocar en una ensaladera. Agregar los granos de maíz dulce. Lavar la manzana y cortar en dados o finas láminas. Añadir a la ensaladera. Preparar el aliño en un cuenco separado con sal, pimienta, orégano seco, vinagre de vino, miel y aceite de oliva. Revolver hasta emulsionar. Volcar sobre las verduras frescas, y revolver para integrar. Servir y terminar con unas pipas de calabaza en la superficie."
$ws.Range("F20").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/ensalada-lechuga-maiz-dulce-manzana-receta-verano-arguinano"

$ws.Range("A21").Value = "17/07/2023 09:23"
$ws.Range("B21").Value = "Recetas fáciles para atravesar la ola de calor en el menú semanal vegetariano del 17 de julio"
$ws.Range("C21").Value = "Inés Vazquez Noya"
$ws.Range("D21").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("F21").Value = "https://www.directoalpaladar.com/recetas-vegetarianas/recetas-faciles-para-atravezar-ola-calor-menu-semanal-vegetariano-17-julio"
